# final update by anjali
$wb = $excel.ActiveWorkbook

# --- Asset sheet: insert a new column (vname/sam) before the old column I ---
$ws = $wb.Worksheets.Item("Asset")

$ws.Columns("I").Insert()
$ws.Range("I1").Value = "vname"
$ws.Range("I2").Value = "sam"

# Column widths for the shifted/new columns (H..N)
$ws.Columns("H:I").ColumnWidth = 12.7109375
$ws.Columns("J").ColumnWidth = 23.5703125
$ws.Columns("K").ColumnWidth = 17.7109375
$ws.Columns("L").ColumnWidth = 18.5703125
$ws.Columns("M").ColumnWidth = 15
$ws.Columns("N").ColumnWidth = 24.28515625

# Hyperlink moved from K2 to L2 (old column K shifted right to L)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("L2"), "mailto:sam@gmail.com") | Out-Null

# This sheet becomes the active / selected tab, scrolled so column F is
# the left-most visible column, with N3 selected.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$win.ScrollRow = 1
$ws.Range("N3").Select()

# --- Sitevisit sheet loses tabSelected (Asset becomes the selected tab) ---
$wsSite = $wb.Worksheets.Item("Sitevisit")
$wsSite.Range("F3").Select()
